# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# mirroring the latest export from the SeniorConnect monitoring feed.

$wb = $excel.ActiveWorkbook

# ---- PIR: rows 97-108 ----
$ws = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-02-01", "20:00:25", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:00:30", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:00:35", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:00:40", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:00:45", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:00:48", "20:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-01", "20:00:55", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:01:00", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:01:05", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:01:08", "20:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-01", "20:01:16", "20:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "20:01:21", "20:00", "Bathroom", "No Motion", "Inactive")
)
$startRow = 97
$endRow = 108

# Every new cell in this block is logged as plain text in the source feed
# (dates, percentages and temperatures alike) - force text formatting on
# the whole block first so Excel's auto-detection doesn't silently turn
# "2026-02-01" into a date serial or "77.5%" into a fraction.
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6)).NumberFormat = "@"

for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $pirRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# ---- Humidity: rows 76-87 ----
$ws = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-02-01", "20:00:27", "20:00", "Bathroom", "77.5%", "Active"),
    @("2026-02-01", "20:00:32", "20:00", "Bathroom", "78.0%", "Active"),
    @("2026-02-01", "20:00:37", "20:00", "Bathroom", "78.1%", "Active"),
    @("2026-02-01", "20:00:42", "20:00", "Bathroom", "77.8%", "Active"),
    @("2026-02-01", "20:00:47", "20:00", "Bathroom", "77.3%", "Active"),
    @("2026-02-01", "20:00:53", "20:00", "Bathroom", "77.9%", "Active"),
    @("2026-02-01", "20:00:58", "20:00", "Bathroom", "77.3%", "Active"),
    @("2026-02-01", "20:01:03", "20:00", "Bathroom", "77.8%", "Active"),
    @("2026-02-01", "20:01:08", "20:00", "Bathroom", "77.1%", "Active"),
    @("2026-02-01", "20:01:13", "20:00", "Bathroom", "77.7%", "Active"),
    @("2026-02-01", "20:01:18", "20:00", "Bathroom", "76.5%", "Active"),
    @("2026-02-01", "20:01:23", "20:00", "Bathroom", "76.0%", "Active")
)
$startRow = 76
$endRow = 87

# Every new cell in this block is logged as plain text in the source feed
# (dates, percentages and temperatures alike) - force text formatting on
# the whole block first so Excel's auto-detection doesn't silently turn
# "2026-02-01" into a date serial or "77.5%" into a fraction.
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6)).NumberFormat = "@"

for ($i = 0; $i -lt $humidityRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $humidityRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# ---- Temperature: rows 76-87 ----
$ws = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-02-01", "20:00:28", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:00:33", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:00:38", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:00:43", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:00:48", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:00:53", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:00:58", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:01:03", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:01:08", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:01:13", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:01:18", "20:00", "Bathroom", "25.1C", "Active"),
    @("2026-02-01", "20:01:23", "20:00", "Bathroom", "25.1C", "Active")
)
$startRow = 76
$endRow = 87

# Every new cell in this block is logged as plain text in the source feed
# (dates, percentages and temperatures alike) - force text formatting on
# the whole block first so Excel's auto-detection doesn't silently turn
# "2026-02-01" into a date serial or "77.5%" into a fraction.
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6)).NumberFormat = "@"

for ($i = 0; $i -lt $temperatureRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $temperatureRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
